# Generate Report for Handback
# Adds a new handback entry (c18ffc4e-765d-44fe-b7e4-3ee8fd7de854) as a new
# row inserted just above the existing last row (ebf5ae0e-...) on each of
# the three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) - columns: A File Name, B Path And Name,
# C Extension, D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Shift the last data row (currently row 3) down to row 4, preserving
# formatting, then populate the freshly-inserted row 3 with the new entry.
$wsOverview.Rows(3).Insert()

$wsOverview.Cells.Item(3, 1).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md"
$wsOverview.Cells.Item(3, 2).Value = "e2e\c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md"
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 5).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 6).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 7).Value = "2016-08-30 08:26:37"

# Resize table + re-create hyperlinks (row insert does not auto-shift them)
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e714b08542f8096a1a81cd0d807b6dba63bd084d/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md", "", "", "e2e\a38377e6-a599-44f8-87da-f903eaf54708.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2/e2e/c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md", "", "", "e2e\c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ad4b108adc31b137963514302ced11a7bb94652/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md", "", "", "e2e\ebf5ae0e-f633-45be-8e2d-22e709d01e40.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - 16 columns, A & I are "Source File Name" /
# "Target File" (both hyperlinked), H & K hold datetime text.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows(3).Insert()

$wsZh.Cells.Item(3, 1).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md"
$wsZh.Cells.Item(3, 2).Value = ".md"
$wsZh.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
$wsZh.Cells.Item(3, 4).Value = "e2e"
$wsZh.Cells.Item(3, 5).Value = "ht"
$wsZh.Cells.Item(3, 6).Value = "True"
$wsZh.Cells.Item(3, 7).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2.zh-cn.xlf"
$wsZh.Cells.Item(3, 8).Value = "2016-08-30 08:26:25"
$wsZh.Cells.Item(3, 9).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md"
$wsZh.Cells.Item(3, 10).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2.zh-cn.xlf"
$wsZh.Cells.Item(3, 11).Value = "2016-08-30 08:27:21"
$wsZh.Cells.Item(3, 12).Value = ""
$wsZh.Cells.Item(3, 13).Value = "True"
$wsZh.Cells.Item(3, 14).Value = ""
$wsZh.Cells.Item(3, 15).Value = "False"
$wsZh.Cells.Item(3, 16).Value = ""

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e714b08542f8096a1a81cd0d807b6dba63bd084d/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md", "", "", "a38377e6-a599-44f8-87da-f903eaf54708.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/072f9a7c1c2450f1ddf9017e90c0be2aeab81f39/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md", "", "", "a38377e6-a599-44f8-87da-f903eaf54708.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2/e2e/c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md", "", "", "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2/e2e/c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md", "", "", "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ad4b108adc31b137963514302ced11a7bb94652/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md", "", "", "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fb7d96d94d013827b97c18952d7225988ddbcffd/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md", "", "", "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md")

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) - same layout as zh-cn.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows(3).Insert()

$wsDe.Cells.Item(3, 1).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md"
$wsDe.Cells.Item(3, 2).Value = ".md"
$wsDe.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
$wsDe.Cells.Item(3, 4).Value = "e2e"
$wsDe.Cells.Item(3, 5).Value = "ht"
$wsDe.Cells.Item(3, 6).Value = "True"
$wsDe.Cells.Item(3, 7).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2.de-de.xlf"
$wsDe.Cells.Item(3, 8).Value = "2016-08-30 08:26:37"
$wsDe.Cells.Item(3, 9).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md"
$wsDe.Cells.Item(3, 10).Value = "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2.de-de.xlf"
$wsDe.Cells.Item(3, 11).Value = "2016-08-30 08:27:41"
$wsDe.Cells.Item(3, 12).Value = ""
$wsDe.Cells.Item(3, 13).Value = "True"
$wsDe.Cells.Item(3, 14).Value = ""
$wsDe.Cells.Item(3, 15).Value = "False"
$wsDe.Cells.Item(3, 16).Value = ""

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e714b08542f8096a1a81cd0d807b6dba63bd084d/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md", "", "", "a38377e6-a599-44f8-87da-f903eaf54708.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fd33a455d4c5c0586039d736b5d3c77f15f683a2/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md", "", "", "a38377e6-a599-44f8-87da-f903eaf54708.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2/e2e/c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md", "", "", "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3ca66c7f4a82cb86a3e91cc3eacf78154cdaf7d2/e2e/c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md", "", "", "c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ad4b108adc31b137963514302ced11a7bb94652/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md", "", "", "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dd32b207923f6f5edf135673210856be55f18d4f/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md", "", "", "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md")

Write-Host "Handback report row added for c18ffc4e-765d-44fe-b7e4-3ee8fd7de854.md"
